# Auto-generated edit script: update 2024 (and a few revised prior-year) crime counts
# per commit "Add data for 2024-11-15".
$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7075
$ws.Range("K3").Value = 7329
$ws.Range("F4").Value = 1921
$ws.Range("H4").Value = 1744
$ws.Range("K4").Value = 1525
$ws.Range("K6").Value = 8072
$ws.Range("F7").Value = 24114
$ws.Range("H7").Value = 26057
$ws.Range("K7").Value = 24522

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 216
$ws.Range("K7").Value = 742
$ws.Range("J8").Value = 1854
$ws.Range("K8").Value = 1598
$ws.Range("K9").Value = 113
$ws.Range("K10").Value = 139
$ws.Range("G11").Value = 364
$ws.Range("K11").Value = 451
$ws.Range("K19").Value = 716
$ws.Range("K22").Value = 77
$ws.Range("K23").Value = 247
$ws.Range("K24").Value = 77
$ws.Range("K27").Value = 230
$ws.Range("K29").Value = 1348
$ws.Range("K30").Value = 96
$ws.Range("K31").Value = 281
$ws.Range("K33").Value = 1049
$ws.Range("K37").Value = 828
$ws.Range("K42").Value = 903
$ws.Range("K43").Value = 200
$ws.Range("F50").Value = 118
$ws.Range("K52").Value = 638
$ws.Range("K54").Value = 478
$ws.Range("K57").Value = 95
$ws.Range("K60").Value = 141
$ws.Range("G63").Value = 294
$ws.Range("H63").Value = 296
$ws.Range("J63").Value = 121
$ws.Range("K63").Value = 70
$ws.Range("K65").Value = 573
$ws.Range("K67").Value = 959
$ws.Range("K71").Value = 76
$ws.Range("K72").Value = 121
$ws.Range("K73").Value = 219
$ws.Range("K76").Value = 330
$ws.Range("K79").Value = 601
$ws.Range("K88").Value = 264
$ws.Range("K89").Value = 364
$ws.Range("K90").Value = 234
$ws.Range("K91").Value = 290
$ws.Range("K96").Value = 265
$ws.Range("K98").Value = 126
$ws.Range("K99").Value = 419
$ws.Range("F101").Value = 24114
$ws.Range("H101").Value = 26057
$ws.Range("K101").Value = 24522

# --- West Ridge ---
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 80
$ws.Range("K7").Value = 265

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 248
$ws.Range("K7").Value = 742

# --- Belmont Cragin ---
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 152
$ws.Range("G4").Value = 27
$ws.Range("G7").Value = 364
$ws.Range("K7").Value = 451

# --- Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 115
$ws.Range("K7").Value = 364

# --- Little Village ---
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 174
$ws.Range("K3").Value = 176
$ws.Range("K6").Value = 233
$ws.Range("K7").Value = 638

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J4").Value = 98
$ws.Range("K6").Value = 533
$ws.Range("J7").Value = 1854
$ws.Range("K7").Value = 1598

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 378
$ws.Range("K4").Value = 52
$ws.Range("K6").Value = 330
$ws.Range("K7").Value = 1049

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K6").Value = 247
$ws.Range("K7").Value = 828

# --- New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 188
$ws.Range("K3").Value = 138
$ws.Range("K6").Value = 210
$ws.Range("K7").Value = 573

# --- Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 174
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 419

# --- Fuller Park ---
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 27
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 96

# --- Gage Park ---
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 281

# --- North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 264
$ws.Range("K4").Value = 54
$ws.Range("K6").Value = 270
$ws.Range("K7").Value = 959

# --- Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K4").Value = 28
$ws.Range("K6").Value = 263
$ws.Range("K7").Value = 478

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 381
$ws.Range("K3").Value = 479
$ws.Range("K4").Value = 62
$ws.Range("K7").Value = 1348

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 216
$ws.Range("K7").Value = 716

# --- River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 63
$ws.Range("K7").Value = 330

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 268
$ws.Range("K6").Value = 336
$ws.Range("K7").Value = 903

# --- Avondale ---
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 139

# --- Dunning ---
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 77

# --- Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 247

# --- Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 75
$ws.Range("K7").Value = 290

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 193
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 601

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 206
$ws.Range("K3").Value = 192

# --- Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 126

# --- Lincoln Square ---
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("F4").Value = 27
$ws.Range("F7").Value = 118

# --- Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 113

# --- Portage Park ---
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 219

# --- Albany Park ---
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 216

# --- United Center ---
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 82
$ws.Range("K7").Value = 264

# --- Edgewater ---
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 230

# --- Washington Heights ---
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 234

# --- Mckinley Park ---
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 95

# --- Morgan Park ---
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 141

# --- Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 200

# --- Clearing ---
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 77

# --- Oakland ---
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 76

# --- Old Town ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 121

Write-Output "Applied 151 cell updates across 40 sheets"